$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-06-18"

# Update the label in A7 ("June (through 06-17)" -> "June (through 06-18)")
$ws.Range("A7").Value = "June (through 06-18)"

# Update the June row (row 7) values for years 2015, 2017-2022 (columns B, D, E, F, G, H, I)
$ws.Range("B7").Value = 10
$ws.Range("D7").Value = 38
$ws.Range("E7").Value = 36
$ws.Range("F7").Value = 30
$ws.Range("G7").Value = 67
$ws.Range("H7").Value = 66
$ws.Range("I7").Value = 92

# Update the Total row (row 8) values for years 2015, 2017-2022 (columns B, D, E, F, G, H, I)
$ws.Range("B8").Value = 118
$ws.Range("D8").Value = 354
$ws.Range("E8").Value = 331
$ws.Range("F8").Value = 234
$ws.Range("G8").Value = 425
$ws.Range("H8").Value = 697
$ws.Range("I8").Value = 755
